$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1: new header labels (replaces the old "B C D" placeholder in A1)
$ws.Range("A1").Value = "TCD ID"
$ws.Range("B1").Value = "TP ID"
$ws.Range("C1").Value = "TC ID"

# Rows 2-3: single-column id values, stored as text (not numbers) -
# force a Text number format before assignment so "1"/"2" are kept
# as literal text rather than being coerced to numeric values.
$ws.Range("A2:A3").NumberFormat = "@"
$ws.Range("A2").Value = "1"
$ws.Range("A3").Value = "2"
